# Auto-generated edit script: applies the exact cell-level changes
# derived from the canonical OOXML diff (rows 2-38 of the single sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 2).Value = 58256
$ws.Cells.Item(4, 1).Value = 130873700
$ws.Cells.Item(4, 2).Value = 57076
$ws.Cells.Item(4, 4).Value = "LC"
$ws.Cells.Item(4, 5).Value = 102613
$ws.Cells.Item(4, 6).Value = "Orre"
$ws.Cells.Item(4, 7).Value = "Lyrurus tetrix"
$ws.Cells.Item(4, 8).Value = "(Linnaeus, 1758)"
$ws.Cells.Item(4, 11).ClearContents()
$ws.Cells.Item(4, 12).ClearContents()
$ws.Cells.Item(4, 13).Value = "färska spår"
$ws.Cells.Item(4, 14).ClearContents()
$ws.Cells.Item(4, 17).Value = 438768
$ws.Cells.Item(4, 18).Value = 6795206
$ws.Cells.Item(4, 50).Value = "Eva Löfqvist, Alfhild Sehlin"
$ws.Cells.Item(5, 1).Value = 130873745
$ws.Cells.Item(5, 2).Value = 79243
$ws.Cells.Item(5, 4).Value = "NT"
$ws.Cells.Item(5, 5).Value = 6425
$ws.Cells.Item(5, 6).Value = "Garnlav"
$ws.Cells.Item(5, 7).Value = "Alectoria sarmentosa"
$ws.Cells.Item(5, 8).Value = "(Ach.) Ach."
$ws.Cells.Item(5, 11).ClearContents()
$ws.Cells.Item(5, 12).ClearContents()
$ws.Cells.Item(5, 13).ClearContents()
$ws.Cells.Item(5, 14).ClearContents()
$ws.Cells.Item(5, 17).Value = 438633
$ws.Cells.Item(5, 18).Value = 6795187
$ws.Cells.Item(5, 50).Value = "Eva Löfqvist"
$ws.Cells.Item(6, 2).Value = 79243
$ws.Cells.Item(7, 2).Value = 79243
$ws.Cells.Item(8, 1).Value = 130873733
$ws.Cells.Item(8, 2).Value = 79243
$ws.Cells.Item(8, 17).Value = 438651
$ws.Cells.Item(8, 18).Value = 6795214
$ws.Cells.Item(9, 1).Value = 130873730
$ws.Cells.Item(9, 2).Value = 79243
$ws.Cells.Item(9, 10).ClearContents()
$ws.Cells.Item(9, 11).ClearContents()
$ws.Cells.Item(9, 14).ClearContents()
$ws.Cells.Item(9, 17).Value = 438606
$ws.Cells.Item(9, 18).Value = 6795190
$ws.Cells.Item(9, 32).ClearContents()
$ws.Cells.Item(10, 1).Value = 130873741
$ws.Cells.Item(10, 2).Value = 79243
$ws.Cells.Item(10, 10).ClearContents()
$ws.Cells.Item(10, 11).ClearContents()
$ws.Cells.Item(10, 14).ClearContents()
$ws.Cells.Item(10, 17).Value = 438767
$ws.Cells.Item(10, 18).Value = 6795135
$ws.Cells.Item(10, 32).ClearContents()
$ws.Cells.Item(11, 2).Value = 57073
$ws.Cells.Item(12, 2).Value = 79243
$ws.Cells.Item(13, 2).Value = 79243
$ws.Cells.Item(14, 2).Value = 79243
$ws.Cells.Item(15, 2).Value = 56456
$ws.Cells.Item(16, 2).Value = 79243
$ws.Cells.Item(17, 2).Value = 57073
$ws.Cells.Item(18, 1).Value = 130873732
$ws.Cells.Item(18, 2).Value = 79243
$ws.Cells.Item(18, 5).Value = 6425
$ws.Cells.Item(18, 6).Value = "Garnlav"
$ws.Cells.Item(18, 7).Value = "Alectoria sarmentosa"
$ws.Cells.Item(18, 8).Value = "(Ach.) Ach."
$ws.Cells.Item(18, 10).ClearContents()
$ws.Cells.Item(18, 11).ClearContents()
$ws.Cells.Item(18, 14).ClearContents()
$ws.Cells.Item(18, 17).Value = 438571
$ws.Cells.Item(18, 18).Value = 6795200
$ws.Cells.Item(18, 32).ClearContents()
$ws.Cells.Item(18, 50).Value = "Eva Löfqvist"
$ws.Cells.Item(19, 1).Value = 130873723
$ws.Cells.Item(19, 2).Value = 79243
$ws.Cells.Item(19, 17).Value = 438879
$ws.Cells.Item(19, 18).Value = 6795213
$ws.Cells.Item(20, 1).Value = 130873719
$ws.Cells.Item(20, 2).Value = 91828
$ws.Cells.Item(20, 5).Value = 5432
$ws.Cells.Item(20, 6).Value = "Granticka"
$ws.Cells.Item(20, 7).Value = "Porodaedalea chrysoloma s.lat."
$ws.Cells.Item(20, 8).ClearContents()
$ws.Cells.Item(20, 10).ClearContents()
$ws.Cells.Item(20, 11).ClearContents()
$ws.Cells.Item(20, 14).ClearContents()
$ws.Cells.Item(20, 17).Value = 438646
$ws.Cells.Item(20, 18).Value = 6795245
$ws.Cells.Item(20, 32).ClearContents()
$ws.Cells.Item(20, 50).Value = "Eva Löfqvist, Alfhild Sehlin"
$ws.Cells.Item(21, 1).Value = 130873699
$ws.Cells.Item(21, 2).Value = 57076
$ws.Cells.Item(21, 4).Value = "LC"
$ws.Cells.Item(21, 5).Value = 102613
$ws.Cells.Item(21, 6).Value = "Orre"
$ws.Cells.Item(21, 7).Value = "Lyrurus tetrix"
$ws.Cells.Item(21, 8).Value = "(Linnaeus, 1758)"
$ws.Cells.Item(21, 11).ClearContents()
$ws.Cells.Item(21, 12).ClearContents()
$ws.Cells.Item(21, 13).Value = "färska spår"
$ws.Cells.Item(21, 14).ClearContents()
$ws.Cells.Item(21, 17).Value = 438808
$ws.Cells.Item(21, 18).Value = 6795184
$ws.Cells.Item(22, 1).Value = 130873726
$ws.Cells.Item(22, 2).Value = 79243
$ws.Cells.Item(22, 4).Value = "NT"
$ws.Cells.Item(22, 5).Value = 6425
$ws.Cells.Item(22, 6).Value = "Garnlav"
$ws.Cells.Item(22, 7).Value = "Alectoria sarmentosa"
$ws.Cells.Item(22, 8).Value = "(Ach.) Ach."
$ws.Cells.Item(22, 11).ClearContents()
$ws.Cells.Item(22, 12).ClearContents()
$ws.Cells.Item(22, 13).ClearContents()
$ws.Cells.Item(22, 14).ClearContents()
$ws.Cells.Item(22, 17).Value = 438662
$ws.Cells.Item(22, 18).Value = 6795157
$ws.Cells.Item(23, 2).Value = 79243
$ws.Cells.Item(24, 2).Value = 79243
$ws.Cells.Item(25, 2).Value = 79243
$ws.Cells.Item(27, 2).Value = 79243
$ws.Cells.Item(28, 2).Value = 57881
$ws.Cells.Item(30, 2).Value = 79243
$ws.Cells.Item(31, 2).Value = 79243
$ws.Cells.Item(32, 2).Value = 57073
$ws.Cells.Item(33, 1).Value = 130873697
$ws.Cells.Item(33, 2).Value = 91829
$ws.Cells.Item(33, 4).Value = "NT"
$ws.Cells.Item(33, 5).Value = 5442
$ws.Cells.Item(33, 6).Value = "Tallticka"
$ws.Cells.Item(33, 7).Value = "Porodaedalea pini"
$ws.Cells.Item(33, 8).Value = "(Brot.) Murrill"
$ws.Cells.Item(33, 12).ClearContents()
$ws.Cells.Item(33, 13).ClearContents()
$ws.Cells.Item(33, 17).Value = 438905
$ws.Cells.Item(33, 18).Value = 6795075
$ws.Cells.Item(33, 50).Value = "Eva Löfqvist, Alfhild Sehlin"
$ws.Cells.Item(34, 2).Value = 79243
$ws.Cells.Item(35, 2).Value = 79243
$ws.Cells.Item(36, 1).Value = 130873703
$ws.Cells.Item(36, 2).Value = 8451
$ws.Cells.Item(36, 4).Value = "LC"
$ws.Cells.Item(36, 5).Value = 106545
$ws.Cells.Item(36, 6).Value = "Mindre märgborre"
$ws.Cells.Item(36, 7).Value = "Tomicus minor"
$ws.Cells.Item(36, 8).Value = "(Hartig, 1834)"
$ws.Cells.Item(36, 12).ClearContents()
$ws.Cells.Item(36, 13).Value = "äldre gnagspår"
$ws.Cells.Item(36, 17).Value = 439003
$ws.Cells.Item(36, 18).Value = 6795150
$ws.Cells.Item(36, 50).Value = "Eva Löfqvist"
$ws.Cells.Item(37, 2).Value = 79243
$ws.Cells.Item(38, 2).Value = 79243
